# The deck originally has 4 slides:
#   1. EY Technology Overview (title slide)
#   2. Elementi di coding (rif. C#, J#, VB.net, Managed C++)
#   3. Data Modeling
#   4. Archiettetture e componenti ERP
#
# The edit keeps only slide 2 ("Elementi di coding ...") and removes the
# other three slides (1, 3, 4). Delete from the highest index down so
# earlier deletions don't shift the indices of slides still to be removed.

$p = $ppt.ActivePresentation

$p.Slides.Item(4).Delete()
$p.Slides.Item(3).Delete()
$p.Slides.Item(1).Delete()
